$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for rows 2-12 is being bumped from 2023-10-08 (45207)
# to 2023-10-09 (45208) as a date serial value.
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = 45208
}
